$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Row 7: Experimental = "false" ---
# Assigning the literal text "false" straight to .Value makes Excel's COM
# layer infer a Boolean (t="b") instead of a shared string, which also
# flips the cell's type away from what the sheet expects. Writing it as a
# quote-prefixed literal in a scratch cell first forces Excel to treat it
# as text, then a values-only paste carries just that text into B7 without
# disturbing B7's existing cell style.
$scratch = $ws.Range("Z1")
$scratch.Value = "'false"
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

# --- Row 8: Date refreshed ---
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# --- Row 17: Description populated ---
$ws.Range("B17").Value = "Categories for interpreting recovery readiness scores from wearable devices"
